$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 117 / 118 hour bumps + new "Extract/Load Image Centerline" note ---
$ws.Range("C117").Value = 5
$ws.Range("C118").Value = 6

$ws.Range("E117:E117").Copy()
$ws.Range("E118").PasteSpecial(-4122)
$ws.Range("E118").Value = "• Extract/Load Image Centerline"

# --- Row 121 relabel: "* Meetings & other" -> "* Patients Database" ---
$oldB121 = "* Meetings & other"
$ws.Range("B121").Value = "* Patients Database"

# --- Insert a fresh data row at 122 (pushes Total/Paid/NotPaid rows down) ---
$ws.Rows.Item(122).Insert()
$ws.Range("E122").Clear()

$ws.Range("B121:C121").Copy()
$ws.Range("B122").PasteSpecial(-4122)
$ws.Range("B122").Value = $oldB121
$ws.Range("C122").Value = 2

# --- Fix up the Total Hours row (now row 123): include the new row in SUM ---
$ws.Range("C123").Formula = "=SUM(C116:C122)"

# --- Fix up Paid/Not Paid rows (the D column here was one row off even
#     before this edit - the insert naturally shifts it by one, but it
#     needs to land two rows down to line up with the Paid/Not Paid labels) ---
$ws.Range("D123").Clear()
$ws.Range("D124").Value = 0
$ws.Range("D125").Formula = "=C123-D124"

$srcFmt = $ws.Range("D123")
$ws.Range("C121:C121").Copy() | Out-Null

# Match the literal "Paid" D-value style (s=3) and the "Not Paid" formula style (s=21)
$ws.Range("D111:D111").Copy()
$ws.Range("D124").PasteSpecial(-4122)
$ws.Range("D112:D112").Copy()
$ws.Range("D125").PasteSpecial(-4122)

$ws.Range("D124").Value = 0
$ws.Range("D125").Formula = "=C123-D124"

# --- Sheet view bookkeeping to match the saved selection/scroll position ---
$ws.Application.ActiveWindow.ScrollRow = 109
$ws.Range("D124").Select()
